# Updated cryptos list values (price + 1h volume change) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to literal text so numeric-looking strings (e.g. "306.40")
# are not auto-coerced into floating point numbers by the Value setter -- this is
# exactly what typing the value into a Text-formatted cell in real Excel would do.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$ws.Range("D2").Value = "42.095.37"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "2.261.99"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "306.40"
$ws.Range("E5").Value = "  -0.27%  "

Set-TextValue $ws.Range("D6") "96.99"
$ws.Range("E6").Value = "  +1.21%  "

Set-TextValue $ws.Range("D7") "0.526"
$ws.Range("E7").Value = "  -0.99%  "

$ws.Range("E8").Value = "  +0.03%  "

Set-TextValue $ws.Range("D9") "0.488"
$ws.Range("E9").Value = "  -1.39%  "

Set-TextValue $ws.Range("D10") "35.04"
$ws.Range("E10").Value = "  -2.36%  "

Set-TextValue $ws.Range("D11") "0.0786"
$ws.Range("E11").Value = "  -2.31%  "

$ws.Range("E12").Value = "  +0.29%  "

Set-TextValue $ws.Range("D13") "6.83"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "2.613.76"
$ws.Range("E14").Value = "  -1.10%  "

Set-TextValue $ws.Range("D15") "14.68"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").Value = "2.256.98"
$ws.Range("E16").Value = "  -1.48%  "

Set-TextValue $ws.Range("D17") "0.788"
$ws.Range("E17").Value = "  -2.04%  "

$ws.Range("D18").Value = "41.975.13"
$ws.Range("E18").Value = "  -0.89%  "

Set-TextValue $ws.Range("D19") "12.22"
$ws.Range("E19").Value = "  -3.41%  "

$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -2.08%  "

Set-TextValue $ws.Range("D21") "5.98"
$ws.Range("E21").Value = "  -0.60%  "

Set-TextValue $ws.Range("D22") "67.56"
$ws.Range("E22").Value = "  -0.66%  "

Set-TextValue $ws.Range("D23") "236.47"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("E24").Value = "  +2.30%  "

Set-TextValue $ws.Range("D25") "2.57"
$ws.Range("E25").Value = "  -1.15%  "

Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.08%  "

Set-TextValue $ws.Range("D27") "23.46"
$ws.Range("E27").Value = "  -2.34%  "

Set-TextValue $ws.Range("D28") "36.95"
$ws.Range("E28").Value = "  +3.03%  "

Set-TextValue $ws.Range("D29") "9.53"
$ws.Range("E29").Value = "  -0.56%  "

Set-TextValue $ws.Range("D30") "2.12"
$ws.Range("E30").Value = "  +1.20%  "

Set-TextValue $ws.Range("D31") "164.35"
$ws.Range("E31").Value = "  +2.15%  "

Set-TextValue $ws.Range("D32") "5.22"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("E33").Value = "  +0.11%  "

Set-TextValue $ws.Range("D35") "17.62"
$ws.Range("E35").Value = "  +2.22%  "

Set-TextValue $ws.Range("D36") "0.0731"
$ws.Range("E36").Value = "  -3.10%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("E38").Value = "  -4.77%  "

$ws.Range("E39").Value = "  -2.25%  "

$ws.Range("E40").Value = "  -1.34%  "

Set-TextValue $ws.Range("D41") "4.09"
$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("E42").Value = "  +1.91%  "

# Rows 43/44 swap rank order: Maker now ranks above EnergySwap.
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.948.92"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "18.88"
$ws.Range("E44").Value = "  -2.79%  "

Set-TextValue $ws.Range("D45") "0.0279"
$ws.Range("E45").Value = "  -2.09%  "

Set-TextValue $ws.Range("D46") "9.97"
$ws.Range("E46").Value = "  -1.92%  "

Set-TextValue $ws.Range("D47") "2.89"
$ws.Range("E47").Value = "  -4.11%  "

Set-TextValue $ws.Range("D48") "53.28"
$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("D49").Value = "2.487.83"
$ws.Range("E49").Value = "  -0.95%  "

Set-TextValue $ws.Range("D50") "71.97"
$ws.Range("E50").Value = "  -1.31%  "

Set-TextValue $ws.Range("D51") "92.14"
$ws.Range("E51").Value = "  -0.11%  "
